$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-10-30 Thursday" "2025-10-31 Friday"

Replace-Text "505÷8=" "235÷8="
Replace-Text "461÷3=" "910÷9="
Replace-Text "124÷4=" "822÷4="
Replace-Text "274÷4=" "432÷2="
Replace-Text "512÷8=" "304÷3="
Replace-Text "377÷6=" "122÷3="
Replace-Text "565÷5=" "940÷8="
Replace-Text "480÷8=" "206÷4="
Replace-Text "933÷4=" "203÷8="
Replace-Text "730÷8=" "158÷2="
Replace-Text "737÷5=" "773÷8="
Replace-Text "279÷3=" "854÷2="
Replace-Text "302÷4=" "724÷3="
Replace-Text "386÷5=" "347÷6="
Replace-Text "435÷2=" "990÷9="
Replace-Text "670÷5=" "209÷2="
Replace-Text "812÷5=" "954÷3="
Replace-Text "156÷5=" "496÷7="
Replace-Text "392÷8=" "204÷9="
Replace-Text "389÷3=" "781÷2="
Replace-Text "571÷8=" "693÷5="
Replace-Text "832÷9=" "445÷6="
Replace-Text "651÷9=" "992÷3="
Replace-Text "234÷4=" "398÷8="
Replace-Text "543÷4=" "295÷3="

Write-Output "Done"
